$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# The importer now writes a trailing blank row (all cells formatted like the
# preceding data row, but empty) instead of failing when the source XLS had
# an empty row. Reproduce that: clone row 3's formatting down into row 4 ...
$ws.Range("A3:L3").Copy() | Out-Null
$ws.Range("A4:L4").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# ... then make sure row 4 carries no real data, except for the single
# whitespace value the import leaves behind in the last column.
$ws.Range("A4:L4").ClearContents() | Out-Null
$ws.Range("L4").Value = " "

# Match the author's final cursor position after the edit.
$ws.Range("L5").Select() | Out-Null
